# Task Sheet updated 07/03/2019
# Applies the changes described by the diff against Task.xlsx:
#  - moves the workbook window down a bit (best effort; yWindow)
#  - renames / inserts several "Inventory*" class names across the
#    classes (sheet2) and methods (sheet3) sheets
#  - updates a couple of dates and adds a new date pair on overAllChart
#    (sheet1), and tweaks a handful of view selections

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Workbook window position (best effort - xWindow/yWindow)
# ---------------------------------------------------------------------
$excel.Windows.Item(1).Top = 630

$wsChart   = $wb.Worksheets.Item("overAllChart")
$wsClasses = $wb.Worksheets.Item("classes")
$wsMethods = $wb.Worksheets.Item("methods")

# ---------------------------------------------------------------------
# 2. overAllChart (sheet1) - date updates
# ---------------------------------------------------------------------
# "Today working" date moved from 1/3/2019 to 7/3/2019
$wsChart.Range("E14").Value = 43649

# New Started/Ended date pair added for the "DataFlow" row
$wsChart.Range("D22").Value = 43558
$wsChart.Range("E22").Value = 43591
$wsChart.Range("D22:E22").NumberFormat = "m/d/yy"

# ---------------------------------------------------------------------
# 3. classes (sheet2) & methods (sheet3) - updated class list
# ---------------------------------------------------------------------
$classList = @(
    "InventoryLoginDetails",
    "InventoryRegisterDetails",
    "InventoryClient",
    "InventoryDdsClient",
    "InventoryDdsClientInfo",
    "InventoryRunningCustomer",
    "InventoryRunningCustomerInfo",
    "InventoryReseller",
    "InventoryResellerInfo",
    "InventoryManuallyInsertingModule",
    "InventoryNormalClient",
    "inventoryCompanyClient",
    "InventoryCustomer",
    "InventoryNormalCustomer",
    "InventoryCompanyCustomer",
    "InventoryEmployee",
    "InventorySalesTeam",
    "InventoryProcurementTeam",
    "InventoryFinanceTeam",
    "InventoryTechnicianTeam",
    "InventoryManagement",
    "InventoryDirector",
    "InventoryManagingDirector",
    "InventoryGeneralManager",
    "InventoryManager",
    "InventoryCatalog",
    "InventoryAddByLaptop",
    "InventoryAddByDesktop",
    "InventoryAddByTablet",
    "InventoryAddByServer",
    "InventoryReport",
    "InventoryTechnicianReportOnCatalog",
    "InventoryClientReport",
    "InventoryCustomerReport",
    "InventoryReportOnTechnicianByManagement",
    "InventoryRunningCustomerReport",
    "InventoryResellerReport",
    "InventoryNeedToBeService",
    "InventoryScrap",
    "InventorySpair",
    "InventoryReadyForSale"
)

# classes sheet: single column G, starting at row 4
for ($i = 0; $i -lt $classList.Count; $i++) {
    $wsClasses.Cells.Item(4 + $i, 7).Value = $classList[$i]
}

# methods sheet: same list replicated across columns A, C and G, starting at row 5
for ($i = 0; $i -lt $classList.Count; $i++) {
    $r = 5 + $i
    $wsMethods.Cells.Item($r, 1).Value = $classList[$i]
    $wsMethods.Cells.Item($r, 3).Value = $classList[$i]
    $wsMethods.Cells.Item($r, 7).Value = $classList[$i]
}

# ---------------------------------------------------------------------
# 4. Sheet view / selection tweaks
# ---------------------------------------------------------------------
$wsClasses.Activate()
$wsClasses.Range("G4:G44").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1

$wsMethods.Activate()
$wsMethods.Range("G5:G45").Select()

$wsChart.Activate()
$wsChart.Range("B14").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
